# Auto commit at 2025-10-30  7:39:50.13
#
# Updates the "Metrics" sheet's monthly/yearly/total figures (column B,
# rows 2-13) to their refreshed values. The "today" sheet pulls several of
# these via direct formulas (=Metrics!Bn) plus running-total formulas in
# columns E/F, so those cascade automatically on recalculation - no need to
# touch them directly. Only the two sheets' cursor/selection positions are
# also updated to match where the editor last clicked.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metrics sheet: refresh the figures
# ---------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 401769.82
$metrics.Range("B3").Value  = 329266.48000000004
$metrics.Range("B4").Value  = 127962.18999999999
$metrics.Range("B5").Value  = 16093
$metrics.Range("B6").Value  = 4768901.29
$metrics.Range("B7").Value  = 4019085.1499999994
$metrics.Range("B8").Value  = 1398564.3299999998
$metrics.Range("B9").Value  = 185094
$metrics.Range("B10").Value = 33234225.090999823
$metrics.Range("B11").Value = 31294306.670000002
$metrics.Range("B12").Value = 11680273.220000004
$metrics.Range("B13").Value = 1282721

# Move the selection/active cell like the editor left it.
[void]$metrics.Activate()
[void]$metrics.Range("F19").Select()

# ---------------------------------------------------------------------
# "today" sheet: just move the selection/active cell.
# All of its B/E/F formulas reference Metrics! (or each other) and its
# A1 "=TODAY()-1" cell is volatile, so they recompute on their own.
# ---------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")
[void]$today.Activate()
[void]$today.Range("G13").Select()
